# ---------------------------------------------------------------------------
# 1) Update "总计" (summary) sheet: insert a new row for 2022-Q3 at the top
#    of the data block, pushing existing quarters down by one row.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$wsSummary = $wb.Worksheets.Item(1)

$wsSummary.Rows.Item(2).Insert()
# Copy formatting (incl. style) of the existing index column cell down into
# the freshly inserted row so A2 keeps the same style as A3:A9.
$wsSummary.Range("A3").Copy()
$wsSummary.Range("A2").PasteSpecial(-4122)
# Remove the bold formatting that Insert() carried over into B2:D2 from row 1.
$wsSummary.Range("B2:D2").ClearFormats()

$wsSummary.Range("A2").Value = 0
$wsSummary.Range("B2").Value = "2022-Q3"
$wsSummary.Range("C2").Value = 18
$wsSummary.Range("D2").Value = 4.9

# ---------------------------------------------------------------------------
# 2) Create the new "2022-Q3" worksheet, positioned right after "总计".
#    We clone an existing quarterly sheet (2021-Q4) so that the header row,
#    column styles (bold header, bordered "#" column) and number formats
#    match the rest of the workbook exactly, then overwrite its data.
# ---------------------------------------------------------------------------
$wsTemplate = $wb.Worksheets.Item("2021-Q4")
$wsTemplate.Copy($null, $wsSummary)
$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# Template sheet has 19 data rows (rows 2-20); 2022-Q3 only needs 18 (rows 2-19).
$wsQ3.Rows.Item(20).Delete()

$q3data = @(
  @('000480', '东方红新动力灵活配置混合', '12.63', '78.30', '7.66', '0.9675', 1),
  @('008269', '大成睿享混合A', '19.80', '66.91', '3.78', '0.7484', 6),
  @('001044', '嘉实新消费股票', '9.52', '82.84', '5.63', '0.5360', 6),
  @('000242', '景顺长城策略精选', '11.62', '90.68', '4.19', '0.4869', 9),
  @('001564', '东方红京东大数据灵活配置混合', '8.84', '73.95', '5.04', '0.4455', 4),
  @('004355', '嘉实丰和灵活配置混合', '8.33', '85.20', '4.73', '0.3940', 7),
  @('003396', '东方红优享红利混合', '14.04', '60.52', '2.44', '0.3426', 6),
  @('169103', '东方红睿轩三年定期开放灵活配置混合', '11.31', '70.03', '2.39', '0.2703', 9),
  @('008704', '广发高股息优享混合A', '2.33', '92.85', '9.06', '0.2111', 1),
  @('004119', '广发创新驱动灵活配置混合', '1.75', '91.75', '9.06', '0.1586', 2),
  @('008270', '大成睿享混合C', '4.02', '66.91', '3.78', '0.1520', 6),
  @('008705', '广发高股息优享混合C', '0.69', '92.85', '9.06', '0.0625', 1),
  @('011698', '南方均衡回报混合A', '4.10', '49.31', '1.49', '0.0611', 10),
  @('460009', '华泰柏瑞量化先行混合A', '4.22', '93.06', '0.98', '0.0414', 9),
  @('014305', '华泰柏瑞中证500指数增强A', '2.20', '34.78', '0.51', '0.0112', 2),
  @('011701', '南方均衡回报混合C', '0.38', '49.31', '1.49', '0.0057', 10),
  @('010246', '华泰柏瑞量化先行混合C', '0.25', '93.06', '0.98', '0.0024', 9),
  @('014306', '华泰柏瑞中证500指数增强C', '0.07', '34.78', '0.51', '0.0004', 2),
)
# Columns D, E, F, G hold numeric-looking figures that are stored as plain
# text in this workbook (e.g. "12.63", "0.9675"); force text format so the
# values aren't reinterpreted as numbers when assigned.
$wsQ3.Range("D2:G19").NumberFormat = "@"
$wsQ3.Range("B2:B19").NumberFormat = "@"

$r = 2
foreach ($row in $q3data) {
    $wsQ3.Cells.Item($r, 2).Value = $row[0]
    $wsQ3.Cells.Item($r, 3).Value = $row[1]
    $wsQ3.Cells.Item($r, 4).Value = $row[2]
    $wsQ3.Cells.Item($r, 5).Value = $row[3]
    $wsQ3.Cells.Item($r, 6).Value = $row[4]
    $wsQ3.Cells.Item($r, 7).Value = $row[5]
    $wsQ3.Cells.Item($r, 8).Value = $row[6]
    $r++
}

# Drop the temporary "@" text format now that the strings are locked in as
# text, so the cells end up unstyled (s attribute) just like the rest of the
# workbook's quarterly sheets.
$wsQ3.Range("B2:G19").ClearFormats()

"2022-Q3 sheet populated"
